$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.719.39'
$ws.Range('E2').Value = '  +5.13%  '
$ws.Range('D3').Value = '2.229.12'
$ws.Range('E3').Value = '  +3.42%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '''228.83'
$ws.Range('E5').Value = '  +0.78%  '
$ws.Range('D6').Value = '''0.624'
$ws.Range('E6').Value = '  -0.54%  '
$ws.Range('D7').Value = '''61.29'
$ws.Range('E7').Value = '  -2.99%  '
$ws.Range('D9').Value = '''0.403'
$ws.Range('E9').Value = '  +2.99%  '
$ws.Range('D10').Value = '''57.98'
$ws.Range('E10').Value = '  -0.66%  '
$ws.Range('E11').Value = '  +4.11%  '
$ws.Range('E12').Value = '  -0.06%  '
$ws.Range('D13').Value = '2.559.68'
$ws.Range('E13').Value = '  +3.37%  '
$ws.Range('D14').Value = '''15.66'
$ws.Range('E14').Value = '  -1.23%  '
$ws.Range('D15').Value = '''21.56'
$ws.Range('E15').Value = '  -1.24%  '
$ws.Range('D16').Value = '''0.794'
$ws.Range('E16').Value = '  -1.06%  '
$ws.Range('E17').Value = '  +1.54%  '
$ws.Range('D18').Value = '2.222.32'
$ws.Range('E18').Value = '  +2.89%  '
$ws.Range('D19').Value = '41.653.09'
$ws.Range('E19').Value = '  +5.24%  '
$ws.Range('D20').Value = '''72.71'
$ws.Range('E20').Value = '  +1.56%  '
$ws.Range('E21').Value = '  +5.45%  '
$ws.Range('D22').Value = '''6.03'
$ws.Range('E22').Value = '  -0.63%  '
$ws.Range('D23').Value = '''248.40'
$ws.Range('E23').Value = '  +8.20%  '
$ws.Range('E24').Value = '  +0.17%  '
$ws.Range('E25').Value = '  -0.14%  '
$ws.Range('D26').Value = '''2.31'
$ws.Range('E26').Value = '  +0.25%  '
$ws.Range('D27').Value = '''9.57'
$ws.Range('E27').Value = '  +1.72%  '
$ws.Range('D28').Value = '''167.84'
$ws.Range('E28').Value = '  -2.49%  '
$ws.Range('D29').Value = '''0.140'
$ws.Range('E29').Value = '  +0.62%  '
$ws.Range('D30').Value = '''19.91'
$ws.Range('E30').Value = '  +0.47%  '
$ws.Range('E31').Value = '  -2.13%  '
$ws.Range('D32').Value = '''2.59'
$ws.Range('E32').Value = '  -3.67%  '
$ws.Range('E33').Value = '  +0.06%  '
$ws.Range('D34').Value = '''5.05'
$ws.Range('E34').Value = '  +8.13%  '
$ws.Range('E35').Value = '  +1.66%  '
$ws.Range('D36').Value = '''0.0624'
$ws.Range('E36').Value = '  +0.89%  '
$ws.Range('B37').Value = 'THORChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D37').Value = '''6.56'
$ws.Range('E37').Value = '  -5.24%  '
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D38').Value = '''3.69'
$ws.Range('E38').Value = '  +2.10%  '
$ws.Range('D39').Value = '''2.37'
$ws.Range('E39').Value = '  -0.89%  '
$ws.Range('D40').Value = '''1.00'
$ws.Range('E40').Value = '  +0.00%  '
$ws.Range('D41').Value = '''0.000236'
$ws.Range('E41').Value = '  +28.58%  '
$ws.Range('E42').Value = '  -3.85%  '
$ws.Range('E43').Value = '  +4.67%  '
$ws.Range('D44').Value = '''8.72'
$ws.Range('E44').Value = '  +12.72%  '
$ws.Range('D45').Value = '''0.0982'
$ws.Range('E45').Value = '  +6.97%  '
$ws.Range('D46').Value = '''99.17'
$ws.Range('E46').Value = '  -3.47%  '
$ws.Range('D47').Value = '1.469.25'
$ws.Range('E47').Value = '  -3.13%  '
$ws.Range('E48').Value = '  -2.36%  '
$ws.Range('D49').Value = '''16.46'
$ws.Range('E49').Value = '  -6.41%  '
$ws.Range('E50').Value = '  -0.87%  '
$ws.Range('E51').Value = '  -0.88%  '
